$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("D2").Value = -0.035
$ws.Range("E2").Value = 0.06367500000000001
$ws.Range("F2").Value = 0.04650000000000001
$ws.Range("G2").Value = 0.08458711825393564
$ws.Range("H2").Value = 0.08458711825393564
$ws.Range("I2").Value = 0.06586801011744306
$ws.Range("J2").Value = 0.04796928613329151
$ws.Range("K2").Value = 3648.94
$ws.Range("L2").Value = 0.02999662953676682
$ws.Range("M2").Value = 1625.856
$ws.Range("N2").Value = 0.04220097854723372
$ws.Range("O2").Value = 0.4455693982362001
$ws.Range("P2").Value = 1619.706
$ws.Range("Q2").Value = 0.04204134816295278
$ws.Range("R2").Value = 0.4438839772646303
$ws.Range("S2").Value = 6.150000000000009
$ws.Range("T2").Value = 0.003782622815304682
$ws.Range("U2").Value = 10357.78
$ws.Range("V2").Value = 0.2688481953979728
$ws.Range("W2").Value = 0.1324582338902148
$ws.Range("X2").Value = 0.06810148051510871
$ws.Range("Y2").Value = 0.06435675337510609
$ws.Range("Z2").Value = 2.070885254104113
$ws.Range("AA2").Value = 0.1242793905125435
$ws.Range("AB2").Value = 0.05718289323928824
$ws.Range("AC2").Value = 0.06426918204249167
$ws.Range("AD2").Value = 31782.8
$ws.Range("AE2").Value = 6.479546318194497
$ws.Range("AF2").Value = 31789.27954631819
$ws.Range("AG2").Value = 21431.49954631819
$ws.Range("AH2").Value = 0.4520931112678351
$ws.Range("AI2").Value = 0.4012076602530734
$ws.Range("AJ2").Value = 0.3574418711178336
$ws.Range("AK2").Value = 0.311159728432836
$ws.Range("AL2").Value = 1333.16
$ws.Range("AM2").Value = 1333.13
$ws.Range("AN2").Value = 3.098206844113877
$ws.Range("AO2").Value = 6.006090791802934
$ws.Range("AP2").Value = 2.08915572492123
$ws.Range("AQ2").Value = 6.006225949457292

# --- Row 3 ---
$ws.Range("D3").Value = 0.08070000000000001
$ws.Range("E3").Value = 0.136
$ws.Range("G3").Value = 0.141280353200883
$ws.Range("H3").Value = 0.141280353200883
$ws.Range("I3").Value = 0.1247240618101545
$ws.Range("J3").Value = 0.0858372062509927
$ws.Range("K3").Value = 6.14
$ws.Range("L3").Value = 0.06777041942604857
$ws.Range("M3").Value = 3.64
$ws.Range("N3").Value = 0.03244206773618538
$ws.Range("O3").Value = 0.5928338762214984
$ws.Range("P3").Value = 3.64
$ws.Range("Q3").Value = 0.03244206773618538
$ws.Range("R3").Value = 0.5928338762214984
$ws.Range("U3").Value = 9.779999999999999
$ws.Range("V3").Value = 0.08716577540106951
$ws.Range("W3").Value = 0.2053511705685619
$ws.Range("X3").Value = 0.06722520694627208
$ws.Range("Y3").Value = 0.1381259636222898
$ws.Range("Z3").Value = 5.367298578199056
$ws.Range("AA3").Value = 0.4607139150675322
$ws.Range("AB3").Value = 0.05705977759072785
$ws.Range("AC3").Value = 0.4036541374768044
$ws.Range("AD3").Value = 41.2
$ws.Range("AF3").Value = 41.2
$ws.Range("AG3").Value = 31.42
$ws.Range("AH3").Value = 0.2685788787483703
$ws.Range("AI3").Value = 0.4051130776794494
$ws.Range("AJ3").Value = 0.2187717588079655
$ws.Range("AK3").Value = 0.3418189730200174
$ws.Range("AL3").Value = 1.2
$ws.Range("AM3").Value = 1.17
$ws.Range("AN3").Value = 3.21875
$ws.Range("AO3").Value = 9.416666666666668
$ws.Range("AP3").Value = 2.4546875
$ws.Range("AQ3").Value = 9.658119658119659

# --- Row 4 ---
$ws.Range("B4").Value = "UnipolSai Assicurazioni S.p.A. (BIT:US)"
$ws.Range("D4").Value = -0.035
$ws.Range("E4").Value = -0.00865
$ws.Range("F4").ClearContents()
$ws.Range("G4").Value = 0.109928359638299
$ws.Range("H4").Value = 0.109928359638299
$ws.Range("I4").Value = 0.08820297488597088
$ws.Range("J4").Value = 0.06368377618976366
$ws.Range("K4").Value = 910.2
$ws.Range("L4").Value = 0.05587579881765779
$ws.Range("M4").Value = 509.2
$ws.Range("N4").Value = 0.06782099094299414
$ws.Range("O4").Value = 0.5594374862667545
$ws.Range("P4").Value = 509.2
$ws.Range("Q4").Value = 0.06782099094299414
$ws.Range("R4").Value = 0.5594374862667545
$ws.Range("S4").Value = 0
$ws.Range("T4").Value = 0
$ws.Range("U4").Value = 718.4
$ws.Range("V4").Value = 0.09568460309003729
$ws.Range("W4").Value = 0.1324582338902148
$ws.Range("X4").Value = 0.06810148051510871
$ws.Range("Y4").Value = 0.06435675337510609
$ws.Range("Z4").Value = 1.952920443101712
$ws.Range("AA4").Value = 0.1243693484149035
$ws.Range("AB4").Value = 0.05718289323928824
$ws.Range("AC4").Value = 0.06718645517561526
$ws.Range("AD4").Value = 2948.1
$ws.Range("AF4").Value = 2948.1
$ws.Range("AG4").Value = 2229.7
$ws.Range("AH4").Value = 0.2819502491368674
$ws.Range("AI4").Value = 0.2878244993995723
$ws.Range("AJ4").Value = 0.2289760415703913
$ws.Range("AK4").Value = 0.2341064435181588
$ws.Range("AL4").Value = 115.2
$ws.Range("AM4").Value = 115.2
$ws.Range("AN4").Value = 1.383824633871573
$ws.Range("AO4").Value = 12.47222222222222
$ws.Range("AP4").Value = 1.046610965076981
$ws.Range("AQ4").Value = 12.47222222222222

# --- Row 5 ---
$ws.Range("B5").Value = "Unipol Gruppo S.p.A. (BIT:UNI)"
$ws.Range("D5").Value = -0.0596
$ws.Range("E5").Value = 0.172
$ws.Range("F5").Value = -0.03
$ws.Range("G5").Value = 0.1385443014995098
$ws.Range("H5").Value = 0.1385443014995098
$ws.Range("I5").Value = 0.09035787333438222
$ws.Range("J5").Value = 0.07328358289962607
$ws.Range("K5").Value = 855.5
$ws.Range("L5").Value = 0.05176879228339405
$ws.Range("M5").Value = 226.906
$ws.Range("N5").Value = 0.06622479059043283
$ws.Range("O5").Value = 0.2652320280537697
$ws.Range("P5").Value = 225.666
$ws.Range("Q5").Value = 0.0658628841607565
$ws.Range("R5").Value = 0.2637825832846289
$ws.Range("S5").Value = 1.240000000000009
$ws.Range("T5").Value = 0.005464818030373851
$ws.Range("U5").Value = 1018.3
$ws.Range("V5").Value = 0.2972010623704871
$ws.Range("W5").Value = 0.123411376062088
$ws.Range("X5").Value = 0.1039150320493368
$ws.Range("Y5").Value = 0.0194963440127512
$ws.Range("Z5").Value = 1.695869464826312
$ws.Range("AA5").Value = 0.1242793905125435
$ws.Range("AB5").Value = 0.06001020847005187
$ws.Range("AC5").Value = 0.06426918204249167
$ws.Range("AD5").Value = 4910.6
$ws.Range("AF5").Value = 4910.6
$ws.Range("AG5").Value = 3892.3
$ws.Range("AH5").Value = 0.5890198994830212
$ws.Range("AI5").Value = 0.3548429054542301
$ws.Range("AJ5").Value = 0.5318366900773372
$ws.Range("AK5").Value = 0.3035997035997036
$ws.Range("AL5").Value = 184.2
$ws.Range("AM5").Value = 184.2
$ws.Range("AN5").Value = 2.288897175351916
$ws.Range("AO5").Value = 8.106406080347449
$ws.Range("AP5").Value = 1.814253752214039
$ws.Range("AQ5").Value = 8.106406080347449

# --- Row 6 ---
$ws.Range("B6").Value = "Assicurazioni Generali S.p.A. (BIT:G)"
$ws.Range("D6").Value = -0.0259
$ws.Range("E6").Value = -0.0274
$ws.Range("F6").Value = 0.04650000000000001
$ws.Range("G6").Value = 0.0696968364612469
$ws.Range("H6").Value = 0.0696968364612469
$ws.Range("I6").Value = 0.05708973927931564
$ws.Range("J6").Value = 0.03705989387268633
$ws.Range("K6").Value = 1862.5
$ws.Range("L6").Value = 0.02099743746709462
$ws.Range("M6").Value = 883.45
$ws.Range("N6").Value = 0.03226424948049245
$ws.Range("O6").Value = 0.4743355704697987
$ws.Range("P6").Value = 881.2
$ws.Range("Q6").Value = 0.03218207781109281
$ws.Range("R6").Value = 0.4731275167785235
$ws.Range("S6").Value = 2.25
$ws.Range("T6").Value = 0.002546833437093214
$ws.Range("U6").Value = 8605.700000000001
$ws.Range("V6").Value = 0.3142865490455304
$ws.Range("W6").Value = 0.06048648999740192
$ws.Range("X6").Value = 0.08459338534066904
$ws.Range("Y6").Value = -0.02410689534326712
$ws.Range("Z6").Value = 2.18687904128203
$ws.Range("AA6").Value = 0.08104550518231408
$ws.Range("AB6").Value = 0.05887532582646126
$ws.Range("AC6").Value = 0.02217017935585282
$ws.Range("AD6").Value = 23865.6
$ws.Range("AE6").Value = 6.479546318194497
$ws.Range("AF6").Value = 23872.07954631819
$ws.Range("AG6").Value = 15266.37954631819
$ws.Range("AH6").Value = 0.4657623253860707
$ws.Range("AI6").Value = 0.4344131464749822
$ws.Range("AJ6").Value = 0.3579617114936688
$ws.Range("AK6").Value = 0.329394613730631
$ws.Range("AL6").Value = 1029.7
$ws.Range("AM6").Value = 1029.7
$ws.Range("AN6").Value = 4.002864748001134
$ws.Range("AO6").Value = 4.912595901718947
$ws.Range("AP6").Value = 2.560557979500312
$ws.Range("AQ6").Value = 4.912595901718947

# --- Row 7 ---
$ws.Range("B7").Value = "Net Insurance S.p.A. (BIT:NET)"
$ws.Range("D7").Value = -0.115
$ws.Range("E7").ClearContents()
$ws.Range("F7").Value = 0.277
$ws.Range("G7").Value = 0.3789473684210526
$ws.Range("H7").Value = 0.3789473684210526
$ws.Range("I7").Value = 0.1915789473684211
$ws.Range("J7").Value = 0.1476881091617934
$ws.Range("K7").Value = 14.6
$ws.Range("L7").Value = 0.3842105263157894
$ws.Range("M7").Value = 2.66
$ws.Range("N7").Value = 0.02706002034587996
$ws.Range("O7").Value = 0.1821917808219178
$ws.Range("P7").Value = -0
$ws.Range("Q7").Value = -0
$ws.Range("R7").Value = -0
$ws.Range("S7").Value = 2.66
$ws.Range("T7").Value = 1
$ws.Range("U7").Value = 5.6
$ws.Range("V7").Value = 0.05696846388606307
$ws.Range("W7").Value = 0.2106782106782107
$ws.Range("X7").Value = 0.06064415318742246
$ws.Range("Y7").Value = 0.1500340574907882
$ws.Range("Z7").Value = 0.4914640455250906
$ws.Range("AA7").Value = 0.07258339560460617
$ws.Range("AB7").Value = 0.055964787645217
$ws.Range("AC7").Value = 0.01661860795938917
$ws.Range("AD7").Value = 17.3
$ws.Range("AF7").Value = 17.3
$ws.Range("AG7").Value = 11.7
$ws.Range("AH7").Value = 0.1496539792387543
$ws.Range("AI7").Value = 0.1759918616480163
$ws.Range("AJ7").Value = 0.1063636363636364
$ws.Range("AK7").Value = 0.1262135922330097
$ws.Range("AL7").Value = 2.86
$ws.Range("AM7").Value = 2.86
$ws.Range("AN7").Value = 2.240932642487047
$ws.Range("AO7").Value = 2.545454545454546
$ws.Range("AP7").Value = 1.515544041450777
$ws.Range("AQ7").Value = 2.545454545454546
